{"js": "// Replace the division-problem answers in the worksheet table.\n// Each old string is unique in the document, so a plain text search +\n// full-run replace is sufficient and keeps existing run formatting intact.\nconst replacements = [\n  [\"649\u00f75=129, 4\", \"559\u00f78=69, 7\"],\n  [\"278\u00f78=34, 6\", \"663\u00f76=110, 3\"],\n  [\"565\u00f78=70, 5\", \"288\u00f79=32, 0\"],\n  [\"755\u00f72=377, 1\", \"682\u00f75=136, 2\"],\n  [\"663\u00f75=132, 3\", \"639\u00f79=71, 0\"],\n  [\"443\u00f76=73, 5\", \"673\u00f77=96, 1\"],\n  [\"897\u00f74=224, 1\", \"534\u00f79=59, 3\"],\n  [\"939\u00f77=134, 1\", \"359\u00f78=44, 7\"],\n  [\"296\u00f74=74, 0\", \"541\u00f74=135, 1\"],\n  [\"257\u00f79=28, 5\", \"342\u00f77=48, 6\"],\n  [\"696\u00f72=348, 0\", \"944\u00f79=104, 8\"],\n  [\"491\u00f77=70, 1\", \"115\u00f73=38, 1\"],\n  [\"994\u00f76=165, 4\", \"289\u00f77=41, 2\"],\n  [\"386\u00f73=128, 2\", \"950\u00f78=118, 6\"],\n  [\"841\u00f73=280, 1\", \"151\u00f78=18, 7\"],\n  [\"408\u00f78=51, 0\", \"434\u00f74=108, 2\"],\n  [\"228\u00f72=114, 0\", \"769\u00f78=96, 1\"],\n  [\"679\u00f76=113, 1\", \"244\u00f74=61, 0\"],\n  [\"252\u00f76=42, 0\", \"930\u00f77=132, 6\"],\n  [\"622\u00f77=88, 6\", \"847\u00f76=141, 1\"],\n  [\"432\u00f79=48, 0\", \"152\u00f76=25, 2\"],\n  [\"541\u00f75=108, 1\", \"869\u00f74=217, 1\"],\n  [\"836\u00f77=119, 3\", \"138\u00f72=69, 0\"],\n  [\"469\u00f79=52, 1\", \"836\u00f74=209, 0\"],\n  [\"571\u00f78=71, 3\", \"995\u00f78=124, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division-problem answers in the worksheet table.\n# Each old string is unique in the document, so a Find/Replace over the\n# whole document content is sufficient and preserves existing run formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"649\u00f75=129, 4\", \"559\u00f78=69, 7\"),\n    @(\"278\u00f78=34, 6\", \"663\u00f76=110, 3\"),\n    @(\"565\u00f78=70, 5\", \"288\u00f79=32, 0\"),\n    @(\"755\u00f72=377, 1\", \"682\u00f75=136, 2\"),\n    @(\"663\u00f75=132, 3\", \"639\u00f79=71, 0\"),\n    @(\"443\u00f76=73, 5\", \"673\u00f77=96, 1\"),\n    @(\"897\u00f74=224, 1\", \"534\u00f79=59, 3\"),\n    @(\"939\u00f77=134, 1\", \"359\u00f78=44, 7\"),\n    @(\"296\u00f74=74, 0\", \"541\u00f74=135, 1\"),\n    @(\"257\u00f79=28, 5\", \"342\u00f77=48, 6\"),\n    @(\"696\u00f72=348, 0\", \"944\u00f79=104, 8\"),\n    @(\"491\u00f77=70, 1\", \"115\u00f73=38, 1\"),\n    @(\"994\u00f76=165, 4\", \"289\u00f77=41, 2\"),\n    @(\"386\u00f73=128, 2\", \"950\u00f78=118, 6\"),\n    @(\"841\u00f73=280, 1\", \"151\u00f78=18, 7\"),\n    @(\"408\u00f78=51, 0\", \"434\u00f74=108, 2\"),\n    @(\"228\u00f72=114, 0\", \"769\u00f78=96, 1\"),\n    @(\"679\u00f76=113, 1\", \"244\u00f74=61, 0\"),\n    @(\"252\u00f76=42, 0\", \"930\u00f77=132, 6\"),\n    @(\"622\u00f77=88, 6\", \"847\u00f76=141, 1\"),\n    @(\"432\u00f79=48, 0\", \"152\u00f76=25, 2\"),\n    @(\"541\u00f75=108, 1\", \"869\u00f74=217, 1\"),\n    @(\"836\u00f77=119, 3\", \"138\u00f72=69, 0\"),\n    @(\"469\u00f79=52, 1\", \"836\u00f74=209, 0\"),\n    @(\"571\u00f78=71, 3\", \"995\u00f78=124, 3\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
